$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Tabla2")

# Row 74 ("Agregar loggins donde hay printstacktrace") is no longer the last
# item (index 26); it shifts up to index 27 because a new item is appended.
$ws.Cells.Item(74, 1).Value = 27

# Add a new row to the table; this extends the table range (and autofilter)
# from A1:E74 to A1:E75, matching the sortState-driven re-sort behaviour.
$newRow = $tbl.ListRows.Add()

# Populate the new last row (row 75) with the new task.
$ws.Cells.Item(75, 1).Value = 26
$ws.Cells.Item(75, 2).Value = "Cambiar los loggin por printstacktrace"
$ws.Cells.Item(75, 3).Value = "no"

# Reflect the new selection (user ended up on the new row, column C).
$ws.Activate()
$ws.Range("C75").Select()
